$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 38: Phụ trách (D38) changes from "Dũng" to "Ngọc Anh"
$ws.Range("D38").Value = "Ngọc Anh"

# Row 41: Phụ trách (D41) changes from "Dũng" to "Ngọc Anh"
$ws.Range("D41").Value = "Ngọc Anh"

# Row 41: Tiến độ (E41) changes from "xong front-end" to new status "xong front-end- Xong Backend"
$ws.Range("E41").Value = "xong front-end- Xong Backend"

# Update the selection shown when the sheet is active
$ws.Range("C48").Select()

$wb.Save()
